$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Clear the event details that were reverted in row 22 (data row 20)
$ws.Range("A22").ClearContents()
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

# Remove the now-unused "Images" column (Q) entirely, shrinking the sheet
# from A1:Q25 to A1:P25
$ws.Range("Q1:Q25").EntireColumn.Delete()
